$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet: "Through 2022-11-01" -> "Through 2022-11-02"
$ws.Name = "Through 2022-11-02"

# Update header label in I1 (shared string): "2022 (through 11-01)" -> "2022 (through 11-02)"
$ws.Range("I1").Value = "2022 (through 11-02)"

# Update data values
$ws.Range("I12").Value = 7
$ws.Range("I14").Value = 1408
